# はじめに: NOR Flash Memoryの呼称を Technical Reference Manualに合わせ修正
#   NOR Flash for Configuration and Code Flash Memory -> Configuration Flash Memory
#   NOR Flash for Data Memory                          -> Data Store Flash Memory
#
# The block-diagram legend boxes on the right-hand side of the slide are
# small rectangles whose text is split across a few centered paragraphs:
#   "NOR" / "Flash Memory"                -> "Configuration" / "Flash Memory"
#   "NOR Flash" / "Memory" / "(32 MByte)" -> "Data Store" / "Flash Memory" / "(32 MByte)"

function Rename-LegendBox($Slide, $ShapeName, $NewText) {
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $shp = $Slide.Shapes.Item($i)
        if ($shp.Name -eq $ShapeName -and $shp.HasTextFrame) {
            $shp.TextFrame.TextRange.Text = $NewText
        }
    }
}

$p = $ppt.ActivePresentation

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    # "NOR" / "Flash Memory"  ->  "Configuration" / "Flash Memory"
    Rename-LegendBox $slide "正方形/長方形 67" "Configuration`rFlash Memory"
    Rename-LegendBox $slide "正方形/長方形 68" "Configuration`rFlash Memory"

    # "NOR Flash" / "Memory" / "(32 MByte)"  ->  "Data Store" / "Flash Memory" / "(32 MByte)"
    Rename-LegendBox $slide "正方形/長方形 69" "Data Store`rFlash Memory`r(32 MByte)"
    Rename-LegendBox $slide "正方形/長方形 70" "Data Store`rFlash Memory`r(32 MByte)"
}
